$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new data row at row 35 (table currently runs rows 7-63,
#    item #29 "LAXEOL PI 5MG  250TAB" is at row 35). The new item
#    "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF." is inserted right
#    before it (alphabetically between IVY PRONT and LAXEOL), pushing
#    every following row down by one.
# ------------------------------------------------------------------

# Copy formatting (fonts/fill/number format/borders) from the row
# above (row 34, "IVY PRONT  SYRUP") onto the freshly inserted row.
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(35).Insert()

# Fix up the row height (a bare insert defaults to a generic height).
$ws.Rows.Item(35).RowHeight = 24.75

# Recreate the per-row merged cell layout used by every item row.
$ws.Range("A35:B35").Merge()
$ws.Range("C35:G35").Merge()
$ws.Range("H35:K35").Merge()
$ws.Range("L35:M35").Merge()
$ws.Range("N35:O35").Merge()

# Populate the new row's values.
$ws.Range("A35").Value = 29
$ws.Range("C35").Value = "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF."
$ws.Range("H35").Value = "3:0"
$ws.Range("L35").Value = 1
$ws.Range("N35").Value = "60.00"
$ws.Range("P35").Value = "12.0000"
$ws.Range("Q35").Value = "0:1"

# ------------------------------------------------------------------
# 2) Column A holds the row's running number (1..58) independent of
#    the item that lands on it. Excel's row-insert shifted the old
#    numbers down along with the rest of each row's content, so
#    restore the simple sequential numbering for every row that moved
#    (rows 36-64, i.e. old rows 35-63 which used to read 29..57).
# ------------------------------------------------------------------
For ($r = 36; $r -le 64; $r++) {
    $ws.Range("A" + $r).Value = $r - 6
}

# ------------------------------------------------------------------
# 3) Update the "سرنجات 3 سم" row (now at row 59): sell price and
#    transactions count changed.
# ------------------------------------------------------------------
$ws.Range("P59").Value = "16.0000"
$ws.Range("Q59").Value = "8:0"

# ------------------------------------------------------------------
# 4) Update the totals row (now row 65) to reflect the new sum of the
#    sell-price column.
# ------------------------------------------------------------------
$ws.Range("N65").Value = 3957.95
